# edit.ps1 -- apply "Minor tweaks to may forecast" changes via Word COM interop
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "The half-way point (50%) of the run is expected be reached on June 17"
#    -> "...is expected on June 17"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "is expected be reached on June 17", $true, $false, $false, $false, $false,
    $true, 1, $false, "is expected on June 17", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "the three percentiles (15%, 25%, 50%)" -> "(15%, 25%, and 50%)"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "(15%, 25%, 50%) and the sea ice coverage", $true, $false, $false, $false, $false,
    $true, 1, $false, "(15%, 25%, and 50%) and the sea ice coverage", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Credits paragraph rework
# ---------------------------------------------------------------------------

# 3a) "... proymundy@gmail.com Jordan Watson jordan.watson@noaa.gov ..."
#     -> "... proymundy@gmail.com, and Jordan Watson jordan.watson@noaa.gov ..."
$credits = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$credits.Find.Execute(
    "Jordan Watson", $true, $false, $false, $false, $false,
    $true, 1, $false, ", and Jordan Watson", 2) | Out-Null

# 3b) Remove "and Zach Liller" (name + mailto hyperlink) entirely -- the
#     credit is dropped from the list.
$zachHyperlink = $null
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    if ($d.Hyperlinks.Item($i).Address -eq "mailto:zachary.liller@alaska.gov") {
        $zachHyperlink = $d.Hyperlinks.Item($i)
        break
    }
}
if ($zachHyperlink -ne $null) {
    $precedingHyperlink = $d.Hyperlinks.Item($zachHyperlink.Index - 1)
    $deleteRange = $d.Range($precedingHyperlink.Range.End, $zachHyperlink.Range.End)
    $deleteRange.Delete()
}

# 3c) "... sean.larson@alaska.gov and Holly Carroll holly.carroll@alaska.gov ..."
#     -> "... sean.larson@alaska.gov, and Holly Carroll holly.carroll@alaska.gov ..."
$credits2 = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$credits2.Find.Execute(
    "and Holly Carroll", $true, $false, $false, $false, $false,
    $true, 1, $false, ", and Holly Carroll", 2) | Out-Null

# 3d) Oxford comma before the final "and" in the financial-support sentence.
$credits3 = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$credits3.Find.Execute(
    "NOAA National Marine Fisheries Service and the Alaska Department",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "NOAA National Marine Fisheries Service, and the Alaska Department", 2) | Out-Null
